$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 288 (shifts existing rows 288-332 down to 289-333)
$ws.Rows.Item(288).Insert()

# Populate the newly inserted row 288 with the new record
$ws.Cells.Item(288, 1).Value = 10
$ws.Cells.Item(288, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(288, 3).Value = "La Araucanía"
$ws.Cells.Item(288, 4).Value = 44504
$ws.Cells.Item(288, 5).Value = 9
$ws.Cells.Item(288, 6).Value = 100112028
$ws.Cells.Item(288, 7).Value = "Sandia"
$ws.Cells.Item(288, 8).Value = "Sin especificar"
$ws.Cells.Item(288, 9).Value = "Primera"
$ws.Cells.Item(288, 10).Value = 550
$ws.Cells.Item(288, 11).Value = 800
$ws.Cells.Item(288, 12).Value = 900
$ws.Cells.Item(288, 13).Value = 864
$ws.Cells.Item(288, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(288, 15).Value = "Perú"
$ws.Cells.Item(288, 16).Value = 864
$ws.Cells.Item(288, 17).Value = 1
$ws.Cells.Item(288, 18).Value = "Hortaliza"
